# Update cryptocurrency price/volume figures per the latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.174.08'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").Value = '1.660.38'
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.83'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5149'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.20%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2641'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.94%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06272'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.77'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07749'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.60%  '

$ws.Range("D12").Value = '1.680.74'
$ws.Range("E12").Value = '  +0.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.440'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.19%  '

$ws.Range("D14").Value = '1.887.77'
$ws.Range("E14").Value = '  -1.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5436'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.37%  '

$ws.Range("D16").Value = '0.0₅8106'
$ws.Range("E16").Value = '  -2.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.83'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.24%  '

$ws.Range("D18").Value = '26.195.15'
$ws.Range("E18").Value = '  -1.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.622'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.78'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.09'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.019'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.64%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.90'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1222'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.227'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.14'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.433'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05977'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.42%  '

$ws.Range("E31").Value = '  -1.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.577'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.257'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.596'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9652'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.421'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.771'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5679'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -7.43%  '

$ws.Range("E39").Value = '  -1.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.975'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8599'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D43").Value = '1.018.55'
$ws.Range("E43").Value = '  -6.84%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.35'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.30%  '

$ws.Range("D45").Value = '1.802.13'
$ws.Range("E45").Value = '  -1.01%  '

$ws.Range("D46").Value = '0.0₈111'
$ws.Range("E46").Value = '  -1.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.68'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.007'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05166'
$ws.Range("D50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.458'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.98%  '

